$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 65, shifting existing rows 65-75 down to 66-76.
$ws.Rows.Item(65).Insert()

# Populate the new row 65 with the new weekly price observation.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 45258
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100101001
$ws.Range("J65").Value = "Arándano (blue)"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 8500
$ws.Range("O65").Value = 8500
$ws.Range("P65").Value = 8500
$ws.Range("Q65").Value = "`$/bandeja 2 kilos"
$ws.Range("R65").Value = "Región del Maule"
$ws.Range("S65").Value = 4250
$ws.Range("T65").Value = 2
